$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row
$ws.Range("A1").Value = "id"
$ws.Range("B1").Value = "category_name"

# Data rows: id, category_name
$categories = @(
    @(1, "Kablovi"),
    @(2, "Kablovski pribor i oprema"),
    @(3, "Kablovski nosači, kanalice"),
    @(4, "Rasveta"),
    @(5, "Led rasveta"),
    @(6, "Prekidači i utičnice"),
    @(7, "Utikaci i razdelnici"),
    @(8, "Osigurači"),
    @(9, "Sklopke"),
    @(10, "Releji"),
    @(11, "Ormani i razvodne kutije"),
    @(12, "Senzori i signalizacija"),
    @(13, "Razno"),
    @(14, "Alati")
)

$row = 2
foreach ($cat in $categories) {
    $ws.Cells.Item($row, 1).Value = $cat[0]
    $ws.Cells.Item($row, 2).Value = $cat[1]
    $row = $row + 1
}

# Widen column B to fit the longer category names
$ws.Columns.Item(2).ColumnWidth = 28.6

# Update selection to match the new data range
$ws.Range("A2:A15").Select()
